$wb = $excel.ActiveWorkbook

# Overview sheet: row 3 (e82bb858 file) status -> "Ready for handoff", handoff date updated
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-03-24 09:18:42"

# zh-cn sheet: row 3 (e82bb858 file) status -> "Ready for handoff", handoff datetime updated
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("E3").Value = "2016-03-24 09:18:38"

# de-de sheet: row 3 (e82bb858 file) status -> "Ready for handoff", handoff datetime updated
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("E3").Value = "2016-03-24 09:18:42"
